$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns at D:E (old D/E -> F/G)
$ws.Range("D1:E1").EntireColumn.Insert()

# New header labels
$ws.Range("D1").Value = "Masse 1 [g]"
$ws.Range("E1").Value = "Auslenkung 1 [cm]"

# New "Masse [g]" values (col B) - symmetric measurement run
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 50
$ws.Range("B4").Value = 100
$ws.Range("B5").Value = 150
$ws.Range("B6").Value = 200
$ws.Range("B7").Value = 250
$ws.Range("B8").Value = 300
$ws.Range("B9").Value = 350
$ws.Range("B10").Value = 400

# New "Auslenkung [cm]" values (col C)
$ws.Range("C2").Value = 19.65
$ws.Range("C3").Value = 18.5
$ws.Range("C4").Value = 17.35
$ws.Range("C5").Value = 16.4
$ws.Range("C6").Value = 15.3
$ws.Range("C7").Value = 14.25
$ws.Range("C8").Value = 13.15
$ws.Range("C9").Value = 12.05
$ws.Range("C10").Value = 11.1

# New "Masse 1 [g]" values (col D) / "Auslenkung 1 [cm]" values (col E)
$ws.Range("D2").Value = 400
$ws.Range("E2").Value = 11.15

$ws.Range("D3").Value = 350
$ws.Range("E3").Value = 12.3

$ws.Range("D4").Value = 300
$ws.Range("E4").Value = 13.35

$ws.Range("D5").Value = 250
$ws.Range("E5").Value = 14.45

$ws.Range("D6").Value = 200
$ws.Range("E6").Value = 15.5

$ws.Range("D7").Value = 150
$ws.Range("E7").Value = 16.65

$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 17.45

$ws.Range("D9").Value = 50
$ws.Range("E9").Value = 18.75

$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 19.8
